# Apply updated crypto market data (price & 1h volume change) to Sheet1
# Rows 46/47 (Stellar / ApeXProtocol) were also reordered/swapped in the source feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.595.87'
$ws.Range('E2').Value = '  +4.74%  '
$ws.Range('D3').Value = '3.628.92'
$ws.Range('E3').Value = '  +4.83%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = "'592.66"
$ws.Range('E5').Value = '  +1.88%  '
$ws.Range('D6').Value = "'192.06"
$ws.Range('E6').Value = '  +3.68%  '
$ws.Range('D7').Value = "'0.645"
$ws.Range('E7').Value = '  +2.06%  '
$ws.Range('D8').Value = '3.620.28'
$ws.Range('E8').Value = '  +4.56%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').Value = "'0.178"
$ws.Range('E10').Value = '  +2.54%  '
$ws.Range('D11').Value = "'0.668"
$ws.Range('E11').Value = '  +3.22%  '
$ws.Range('D12').Value = "'58.41"
$ws.Range('E12').Value = '  +3.22%  '
$ws.Range('D13').Value = "'0.0000290"
$ws.Range('E13').Value = '  +3.86%  '
$ws.Range('D14').Value = "'9.93"
$ws.Range('E14').Value = '  +5.08%  '
$ws.Range('D15').Value = '4.214.17'
$ws.Range('E15').Value = '  +5.41%  '
$ws.Range('D16').Value = "'19.78"
$ws.Range('E16').Value = '  +5.82%  '
$ws.Range('D17').Value = '3.626.25'
$ws.Range('E17').Value = '  +5.06%  '
$ws.Range('D18').Value = '70.554.36'
$ws.Range('E18').Value = '  +4.79%  '
$ws.Range('D19').Value = "'12.69"
$ws.Range('E19').Value = '  +4.59%  '
$ws.Range('E20').Value = '  +0.52%  '
$ws.Range('E21').Value = '  +4.13%  '
$ws.Range('D22').Value = "'488.45"
$ws.Range('E22').Value = '  +1.48%  '
$ws.Range('D23').Value = "'19.57"
$ws.Range('E23').Value = '  +16.60%  '
$ws.Range('D24').Value = "'5.35"
$ws.Range('E24').Value = '  -1.62%  '
$ws.Range('E25').Value = '  +0.71%  '
$ws.Range('D26').Value = "'91.05"
$ws.Range('E26').Value = '  +1.47%  '
$ws.Range('D27').Value = "'3.15"
$ws.Range('E27').Value = '  +6.92%  '
$ws.Range('D28').Value = "'11.36"
$ws.Range('E28').Value = '  +2.88%  '
$ws.Range('D29').Value = "'9.70"
$ws.Range('E29').Value = '  +5.68%  '
$ws.Range('D30').Value = "'33.15"
$ws.Range('E30').Value = '  +5.45%  '
$ws.Range('D31').Value = "'7.88"
$ws.Range('E31').Value = '  +10.30%  '
$ws.Range('D32').Value = "'629.20"
$ws.Range('E32').Value = '  +6.31%  '
$ws.Range('E33').Value = '  +5.10%  '
$ws.Range('E34').Value = '  +7.22%  '
$ws.Range('D35').Value = "'66.11"
$ws.Range('E35').Value = '  +2.70%  '
$ws.Range('E36').Value = '  +7.33%  '
$ws.Range('D37').Value = "'39.03"
$ws.Range('E37').Value = '  +6.65%  '
$ws.Range('D38').Value = '0.0₃0813'
$ws.Range('E38').Value = '  +5.19%  '
$ws.Range('E39').Value = '  -1.03%  '
$ws.Range('D40').Value = "'0.999"
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('D41').Value = "'3.60"
$ws.Range('E41').Value = '  +0.80%  '
$ws.Range('D42').Value = '3.306.18'
$ws.Range('E42').Value = '  +3.16%  '
$ws.Range('E43').Value = '  +6.84%  '
$ws.Range('D44').Value = "'2.79"
$ws.Range('E44').Value = '  +9.78%  '
$ws.Range('D45').Value = "'0.0450"
$ws.Range('E45').Value = '  +4.65%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = "'3.31"
$ws.Range('E46').Value = '  +3.11%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').Value = "'0.139"
$ws.Range('E47').Value = '  +2.65%  '
$ws.Range('E48').Value = '  +3.97%  '
$ws.Range('D49').Value = "'2.74"
$ws.Range('E49').Value = '  -1.96%  '
$ws.Range('D50').Value = "'3.30"
$ws.Range('E50').Value = '  +3.45%  '
$ws.Range('D51').Value = "'1.00"
$ws.Range('E51').Value = '  +0.23%  '
